$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 28: extend the time range in B28, keep C28 text, add new D28 text,
# and bump hours in G28 from 2 to 3 (and H3's SUM formula recalculates automatically).
$ws.Range("D28").Value = "Ihan hyvä meininki, tässsä se ymmärrys karttuu kun korjaa virheitä."
$ws.Range("D28").WrapText = $true
$ws.Range("B28").Value = "14.00-16.00, 17.45-18.45, 19.00-"
$ws.Range("G28").Value = 3

# The added D28 text wraps onto two lines, so the row grows taller.
$ws.Rows.Item(28).RowHeight = 29

# Adjust the view to reflect where the user was working.
$ws.Range("B28").Select()
